$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Neutrophils"-target rows (original rows 7 and 4).
# Delete bottom-up so row numbers for the remaining rows stay stable.
$ws.Rows(7).Delete()
$ws.Rows(4).Delete()

# Rename the sending cluster "Resolving-Mac" to "MuSCs" for the rows that
# used it (now rows 4 and 5 after the deletions above).
$ws.Range("A4").Value = "MuSCs"
$ws.Range("A5").Value = "MuSCs"

# Row 3's target cluster becomes "MuSCs" (was already MuSCs, index just shifts
# in the shared-string table once "Resolving-Mac" is renamed) -- no value
# change needed, but row 5's target cluster becomes "MuSCs" too.
$ws.Range("D5").Value = "MuSCs"

# Updated TPM-derived numeric values for row 2 (FAPs -> Cxcl13/Ackr4 -> ECs)
$ws.Range("G2").Value = 7.767740666666666
$ws.Range("I2").Value = 0.957755623847744
$ws.Range("J2").Value = 0.9714348434930534
$ws.Range("M2").Value = 0.1631145
$ws.Range("N2").Value = 0.326229
$ws.Range("O2").Value = 0.7212828052797984
$ws.Range("P2").Value = 0.7212828052797984
$ws.Range("Q2").Value = 1.267031134973
$ws.Range("R2").Value = 7.602186809837999
$ws.Range("S2").Value = 0.6908126631414042
$ws.Range("T2").Value = 0.7006792490612115

# Updated TPM-derived numeric values for row 3 (FAPs -> Cxcl13/Ackr4 -> MuSCs)
$ws.Range("G3").Value = 7.767740666666666
$ws.Range("I3").Value = 0.957755623847744
$ws.Range("J3").Value = 0.9714348434930534
$ws.Range("M3").Value = 0.0630305
$ws.Range("O3").Value = 0.2787171947202017
$ws.Range("P3").Value = 0.2787171947202017
$ws.Range("S3").Value = 0.2669429607063399
$ws.Range("T3").Value = 0.270755594431842

# Updated TPM-derived numeric values for row 4 (MuSCs -> Cxcl13/Ackr4 -> ECs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.342617
$ws.Range("H4").Value = 0.685234
$ws.Range("I4").Value = 0.04224437615225601
$ws.Range("J4").Value = 0.02856515650694651
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.1631145
$ws.Range("N4").Value = 0.326229
$ws.Range("O4").Value = 0.7212828052797984
$ws.Range("P4").Value = 0.7212828052797984
$ws.Range("Q4").Value = 0.0558858006465
$ws.Range("R4").Value = 0.223543202586
$ws.Range("S4").Value = 0.03047014213839423
$ws.Range("T4").Value = 0.02060355621858687

# Updated TPM-derived numeric values for row 5 (MuSCs -> Cxcl13/Ackr4 -> MuSCs)
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.342617
$ws.Range("H5").Value = 0.685234
$ws.Range("I5").Value = 0.04224437615225601
$ws.Range("J5").Value = 0.02856515650694651
$ws.Range("M5").Value = 0.0630305
$ws.Range("O5").Value = 0.2787171947202017
$ws.Range("P5").Value = 0.2787171947202017
$ws.Range("Q5").Value = 0.0215953208185
$ws.Range("R5").Value = 0.08638128327400001
$ws.Range("S5").Value = 0.01177423401386178
$ws.Range("T5").Value = 0.007961600288359646
